$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new reservation row (row 2) mirroring the existing row 1 layout:
# Cliente (email), Alojamiento, FechaInicio, FechaFin, PrecioTotal
# FechaInicio/FechaFin/PrecioTotal must stay as plain text (like row 1),
# not get auto-converted to dates/numbers, so force text format first and
# restore the default "Normal" style afterwards (matches row 1 styling).

$ws.Range("A2").Value = "juliddv1@gmail.com"
$ws.Range("B2").Value = "Finca Bella Vista"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "2025-02-02"
$ws.Range("C2").Style = "Normal"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2025-02-05"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2475.0"
$ws.Range("E2").Style = "Normal"
